$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 987.5
$ws.Range("I19").Value = 805.55554
$ws.Range("J19").Value = 1533.3334
$ws.Range("K19").Value = 805.55554
$ws.Range("L19").Value = 1533.3334
$ws.Range("M19").Value = -630.55554
$ws.Range("N19").Value = -1883.3334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 18646.295
$ws.Range("J69").Value = 18646.295
$ws.Range("L69").Value = 55938.88499999999
$ws.Range("N69").Value = -57686.88499999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 18646.295
$ws.Range("J72").Value = 18646.295
$ws.Range("L72").Value = 167816.655
$ws.Range("N72").Value = -176552.655

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 4364.5
$ws.Range("J88").Value = 2729.3333
$ws.Range("L88").Value = 2729.3333
$ws.Range("N88").Value = -3541.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 4364.5
$ws.Range("J91").Value = 2729.3333
$ws.Range("L91").Value = 2729.3333
$ws.Range("N91").Value = -5537.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5464.788
$ws.Range("I138").Value = 4222.067
$ws.Range("J138").Value = 6500.3887
$ws.Range("K138").Value = 12666.201
$ws.Range("L138").Value = 19501.1661
$ws.Range("M138").Value = -7526.201000000001
$ws.Range("N138").Value = -29781.1661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 33164.43
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 33164.43
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 33164.43
$ws.Range("M32").ClearContents() | Out-Null
$ws.Range("N32").Value = -33738.43

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 7059.143
$ws.Range("I74").Value = 2683
$ws.Range("J74").Value = 17999.5
$ws.Range("K74").Value = 2683
$ws.Range("L74").Value = 17999.5
$ws.Range("M74").Value = -1809
$ws.Range("N74").Value = -19747.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 183291.5
$ws.Range("J76").Value = 183291.5
$ws.Range("L76").Value = 183291.5
$ws.Range("N76").Value = -183967.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 7059.143
$ws.Range("I77").Value = 2683
$ws.Range("J77").Value = 17999.5
$ws.Range("K77").Value = 13415
$ws.Range("L77").Value = 89997.5
$ws.Range("M77").Value = -9047
$ws.Range("N77").Value = -98733.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H79").Value = 183291.5
$ws.Range("J79").Value = 183291.5
$ws.Range("L79").Value = 183291.5
$ws.Range("N79").Value = -185631.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 465
$ws.Range("I102").Value = 465
$ws.Range("K102").Value = 465
$ws.Range("M102").Value = 1157

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2681.1052
$ws.Range("I132").Value = 1597.3846
$ws.Range("K132").Value = 4792.1538
$ws.Range("M132").Value = -2262.1538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 128000
$ws.Range("J23").Value = 128000
$ws.Range("L23").Value = 128000
$ws.Range("N23").Value = -128566

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1605.0769
$ws.Range("I86").Value = 911
$ws.Range("J86").Value = 2200
$ws.Range("K86").Value = 911
$ws.Range("L86").Value = 2200
$ws.Range("M86").Value = 212
$ws.Range("N86").Value = -4446

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1605.0769
$ws.Range("I89").Value = 911
$ws.Range("J89").Value = 2200
$ws.Range("K89").Value = 4555
$ws.Range("L89").Value = 11000
$ws.Range("M89").Value = 1061
$ws.Range("N89").Value = -22232

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 835.4783
$ws.Range("I94").Value = 548.3158
$ws.Range("J94").Value = 2199.5
$ws.Range("K94").Value = 548.3158
$ws.Range("L94").Value = 2199.5
$ws.Range("M94").Value = -97.31579999999997
$ws.Range("N94").Value = -3101.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3355.4736
$ws.Range("I99").Value = 2183.6
$ws.Range("K99").Value = 2183.6
$ws.Range("M99").Value = -685.5999999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10453.6875
$ws.Range("I31").Value = 4635.3076
$ws.Range("K31").Value = 4635.3076
$ws.Range("M31").Value = -4340.3076

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 10453.6875
$ws.Range("I34").Value = 4635.3076
$ws.Range("K34").Value = 4635.3076
$ws.Range("M34").Value = -4433.3076

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 198291.5
$ws.Range("J81").Value = 198291.5
$ws.Range("L81").Value = 198291.5
$ws.Range("N81").Value = -200287.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H84").Value = 198291.5
$ws.Range("J84").Value = 198291.5
$ws.Range("L84").Value = 594874.5
$ws.Range("N84").Value = -604858.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4144.115
$ws.Range("I132").Value = 2234.6428
$ws.Range("J132").Value = 6371.8335
$ws.Range("K132").Value = 6703.928400000001
$ws.Range("L132").Value = 19115.5005
$ws.Range("M132").Value = -4173.928400000001
$ws.Range("N132").Value = -24175.5005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 722.5417
$ws.Range("I2").Value = 953.8333
$ws.Range("K2").Value = 953.8333
$ws.Range("M2").Value = -840.8333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6926.6665
$ws.Range("I80").Value = 6273.8335
$ws.Range("J80").Value = 7448.933
$ws.Range("K80").Value = 6273.8335
$ws.Range("L80").Value = 7448.933
$ws.Range("M80").Value = -5275.8335
$ws.Range("N80").Value = -9444.933000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 6926.6665
$ws.Range("I83").Value = 6273.8335
$ws.Range("J83").Value = 7448.933
$ws.Range("K83").Value = 31369.1675
$ws.Range("L83").Value = 37244.665
$ws.Range("M83").Value = -26377.1675
$ws.Range("N83").Value = -47228.665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 296.22223
$ws.Range("I97").Value = 303.3
$ws.Range("J97").Value = 287.375
$ws.Range("K97").Value = 303.3
$ws.Range("L97").Value = 287.375
$ws.Range("M97").Value = 192.7
$ws.Range("N97").Value = -1279.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3236.353
$ws.Range("I102").Value = 1718.091
$ws.Range("J102").Value = 6019.8335
$ws.Range("K102").Value = 1718.091
$ws.Range("L102").Value = 6019.8335
$ws.Range("M102").Value = -96.09099999999989
$ws.Range("N102").Value = -9263.833500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6660.0884
$ws.Range("I132").Value = 6017.8184
$ws.Range("K132").Value = 18053.4552
$ws.Range("M132").Value = -15523.4552

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7880.5454
$ws.Range("I7").Value = 7782.5
$ws.Range("K7").Value = 7782.5
$ws.Range("M7").Value = -7670.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5412.7144
$ws.Range("I40").Value = 4314.8335
$ws.Range("K40").Value = 4314.8335
$ws.Range("M40").Value = -4178.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 7880.5454
$ws.Range("I126").Value = 7782.5
$ws.Range("K126").Value = 23347.5
$ws.Range("M126").Value = -20877.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5676.846
$ws.Range("I122").Value = 5083.2
$ws.Range("K122").Value = 15249.6
$ws.Range("M122").Value = -12799.6
